$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  33"
$ws.Range("C9").Value = "Report Covering the Week  8/14/2023  Through  8/20/2023"

# --- Numeric -> text placeholder conversions (style 15/16 -> style 14 "General" text) ---
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "***.*"
$ws.Range("E14").NumberFormat = "General"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "0"
$ws.Range("C26").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("F26").NumberFormat = "@"
$ws.Range("F26").Value = "0"
$ws.Range("F26").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("E27").NumberFormat = "General"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "***.*"
$ws.Range("E28").NumberFormat = "General"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "***.*"
$ws.Range("E29").NumberFormat = "General"

# --- Text -> numeric conversions (style 14 -> style 15 "#,##0") ---
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"

# --- Pure numeric value updates (style unchanged) ---
$ws.Range("N14").Value = -84.615384615384
$ws.Range("M15").Value = -21.428571428571
$ws.Range("N15").Value = -57.692307692307
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 133.333333333333
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 42.857142857142
$ws.Range("I16").Value = 111
$ws.Range("J16").Value = 119
$ws.Range("K16").Value = -6.72268907563
$ws.Range("L16").Value = 30.588235294117
$ws.Range("M16").Value = -47.887323943662
$ws.Range("N16").Value = -83.027522935779
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -20
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 29
$ws.Range("H17").Value = -3.448275862068
$ws.Range("I17").Value = 281
$ws.Range("J17").Value = 287
$ws.Range("K17").Value = -2.090592334494
$ws.Range("L17").Value = 6.844106463878
$ws.Range("M17").Value = 46.354166666666
$ws.Range("N17").Value = 12.4
$ws.Range("C18").Value = 4
$ws.Range("E18").Value = -42.857142857142
$ws.Range("G18").Value = 28
$ws.Range("H18").Value = -21.428571428571
$ws.Range("I18").Value = 138
$ws.Range("J18").Value = 113
$ws.Range("K18").Value = 22.123893805309
$ws.Range("L18").Value = 56.818181818181
$ws.Range("M18").Value = -39.473684210526
$ws.Range("N18").Value = -86.227544910179
$ws.Range("C19").Value = 14
$ws.Range("E19").Value = -6.666666666666
$ws.Range("F19").Value = 77
$ws.Range("G19").Value = 59
$ws.Range("H19").Value = 30.508474576271
$ws.Range("I19").Value = 432
$ws.Range("J19").Value = 418
$ws.Range("K19").Value = 3.34928229665
$ws.Range("L19").Value = 52.650176678445
$ws.Range("M19").Value = 55.395683453237
$ws.Range("N19").Value = 14.893617021276
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 3
$ws.Range("E20").Value = 233.333333333333
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 22
$ws.Range("H20").Value = 40.90909090909
$ws.Range("I20").Value = 206
$ws.Range("J20").Value = 151
$ws.Range("K20").Value = 36.423841059602
$ws.Range("L20").Value = 119.148936170213
$ws.Range("M20").Value = -14.876033057851
$ws.Range("N20").Value = -90.480591497227
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 13.157894736842
$ws.Range("G21").Value = 153
$ws.Range("H21").Value = 16.339869281045
$ws.Range("I21").Value = 1181
$ws.Range("J21").Value = 1099
$ws.Range("K21").Value = 7.461328480436
$ws.Range("L21").Value = 39.928909952606
$ws.Range("M21").Value = 0.425170068027
$ws.Range("N21").Value = -73.667781493868
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 46
$ws.Range("E24").Value = -21.739130434782
$ws.Range("F24").Value = 112
$ws.Range("G24").Value = 134
$ws.Range("H24").Value = -16.417910447761
$ws.Range("I24").Value = 793
$ws.Range("J24").Value = 952
$ws.Range("K24").Value = -16.701680672268
$ws.Range("L24").Value = 31.946755407653
$ws.Range("M24").Value = 47.124304267161
$ws.Range("C25").Value = 10
$ws.Range("D25").Value = 8
$ws.Range("F25").Value = 51
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = 18.60465116279
$ws.Range("I25").Value = 405
$ws.Range("J25").Value = 421
$ws.Range("K25").Value = -3.800475059382
$ws.Range("L25").Value = 40.625
$ws.Range("M25").Value = -1.459854014598
$ws.Range("H26").Value = -100
$ws.Range("I26").Value = 21
$ws.Range("K26").Value = 5
$ws.Range("L26").Value = -40
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = -37.5
$ws.Range("I27").Value = 32
$ws.Range("K27").Value = -13.513513513513
$ws.Range("L27").Value = -3.030303030303
$ws.Range("M28").Value = -76.666666666666
$ws.Range("N28").Value = -82.926829268292
$ws.Range("M29").Value = -71.428571428571
$ws.Range("N29").Value = -85
